$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 188
$ws.Range("D188").Value = 44637
$ws.Range("D188").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("K188").Value = 'August Red'
$ws.Range("L188").Value = 'Especial'
$ws.Range("M188").Value = 20
$ws.Range("N188").Value = 450000
$ws.Range("O188").Value = 460000
$ws.Range("P188").Value = 455000
$ws.Range("Q188").Value = '$/bins (420 kilos)'
$ws.Range("S188").Value = 1083
$ws.Range("T188").Value = 420

# Row 189
$ws.Range("D189").Value = 44637
$ws.Range("D189").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("K189").Value = 'August Red'
$ws.Range("L189").Value = 'Primera'
$ws.Range("M189").Value = 20
$ws.Range("N189").Value = 420000
$ws.Range("O189").Value = 430000
$ws.Range("P189").Value = 425000
$ws.Range("Q189").Value = '$/bins (420 kilos)'
$ws.Range("S189").Value = 1012
$ws.Range("T189").Value = 420

# Row 190
$ws.Range("D190").Value = 44272
$ws.Range("D190").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("K190").Value = 'June Pearl'
$ws.Range("L190").Value = 'Especial'
$ws.Range("M190").Value = 160
$ws.Range("N190").Value = 21500
$ws.Range("O190").Value = 22000
$ws.Range("P190").Value = 21750
$ws.Range("Q190").Value = '$/caja 18 kilos empedrada'
$ws.Range("S190").Value = 1208
$ws.Range("T190").Value = 18

# Row 191
$ws.Range("D191").Value = 44272
$ws.Range("D191").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("K191").Value = 'June Pearl'
$ws.Range("L191").Value = 'Primera'
$ws.Range("M191").Value = 240
$ws.Range("N191").Value = 19500
$ws.Range("O191").Value = 20000
$ws.Range("P191").Value = 19750
$ws.Range("Q191").Value = '$/caja 18 kilos empedrada'
$ws.Range("S191").Value = 1097
$ws.Range("T191").Value = 18

# Row 192
$ws.Range("D192").Value = 44615
$ws.Range("D192").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("K192").Value = 'August Red'
$ws.Range("L192").Value = 'Primera'
$ws.Range("M192").Value = 20
$ws.Range("N192").Value = 320000
$ws.Range("O192").Value = 330000
$ws.Range("P192").Value = 325000
$ws.Range("Q192").Value = '$/bins (420 kilos)'
$ws.Range("S192").Value = 774
$ws.Range("T192").Value = 420

# Row 193
$ws.Range("D193").Value = 44615
$ws.Range("D193").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("K193").Value = 'August Red'
$ws.Range("L193").Value = 'Segunda'
$ws.Range("M193").Value = 20
$ws.Range("N193").Value = 270000
$ws.Range("O193").Value = 280000
$ws.Range("P193").Value = 275000
$ws.Range("Q193").Value = '$/bins (420 kilos)'
$ws.Range("S193").Value = 655
$ws.Range("T193").Value = 420

# Row 194
$ws.Range("D194").Value = 44615
$ws.Range("D194").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("K194").Value = 'June Pearl'
$ws.Range("L194").Value = 'Especial'
$ws.Range("M194").Value = 16
$ws.Range("N194").Value = 360000
$ws.Range("O194").Value = 370000
$ws.Range("P194").Value = 365000
$ws.Range("Q194").Value = '$/bins (420 kilos)'
$ws.Range("S194").Value = 869
$ws.Range("T194").Value = 420

# Row 195
$ws.Range("D195").Value = 44615
$ws.Range("D195").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("K195").Value = 'June Pearl'
$ws.Range("L195").Value = 'Primera'
$ws.Range("M195").Value = 20
$ws.Range("N195").Value = 330000
$ws.Range("O195").Value = 340000
$ws.Range("P195").Value = 335000
$ws.Range("Q195").Value = '$/bins (420 kilos)'
$ws.Range("S195").Value = 798
$ws.Range("T195").Value = 420

# Row 196
$ws.Range("D196").Value = 44615
$ws.Range("D196").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("K196").Value = 'Venus'
$ws.Range("L196").Value = 'Especial'
$ws.Range("M196").Value = 28
$ws.Range("N196").Value = 320000
$ws.Range("O196").Value = 330000
$ws.Range("P196").Value = 323571
$ws.Range("Q196").Value = '$/bins (420 kilos)'
$ws.Range("S196").Value = 770
$ws.Range("T196").Value = 420

# Row 197
$ws.Range("D197").Value = 44615
$ws.Range("D197").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("K197").Value = 'Venus'
$ws.Range("L197").Value = 'Primera'
$ws.Range("M197").Value = 20
$ws.Range("N197").Value = 300000
$ws.Range("O197").Value = 310000
$ws.Range("P197").Value = 305000
$ws.Range("Q197").Value = '$/bins (420 kilos)'
$ws.Range("S197").Value = 726
$ws.Range("T197").Value = 420

# Row 198
$ws.Range("D198").Value = 44258
$ws.Range("D198").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("K198").Value = 'August Red'
$ws.Range("L198").Value = 'Especial'
$ws.Range("M198").Value = 240
$ws.Range("N198").Value = 19500
$ws.Range("O198").Value = 20000
$ws.Range("P198").Value = 19750
$ws.Range("Q198").Value = '$/caja 16 kilos empedrada'
$ws.Range("S198").Value = 1234
$ws.Range("T198").Value = 16

# Row 199
$ws.Range("D199").Value = 44258
$ws.Range("D199").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("K199").Value = 'August Red'
$ws.Range("L199").Value = 'Primera'
$ws.Range("M199").Value = 300
$ws.Range("N199").Value = 17500
$ws.Range("O199").Value = 18000
$ws.Range("P199").Value = 17750
$ws.Range("Q199").Value = '$/caja 16 kilos empedrada'
$ws.Range("S199").Value = 1109
$ws.Range("T199").Value = 16

# Row 200
$ws.Range("D200").Value = 44258
$ws.Range("D200").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("K200").Value = 'August Red'
$ws.Range("L200").Value = 'Segunda'
$ws.Range("M200").Value = 180
$ws.Range("N200").Value = 15500
$ws.Range("O200").Value = 16000
$ws.Range("P200").Value = 15750
$ws.Range("Q200").Value = '$/caja 16 kilos empedrada'
$ws.Range("S200").Value = 984
$ws.Range("T200").Value = 16

# Row 201
$ws.Range("D201").Value = 44258
$ws.Range("D201").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("K201").Value = 'June Pearl'
$ws.Range("L201").Value = 'Especial'
$ws.Range("M201").Value = 20
$ws.Range("N201").Value = 390000
$ws.Range("O201").Value = 400000
$ws.Range("P201").Value = 395000
$ws.Range("Q201").Value = '$/bins (420 kilos)'
$ws.Range("S201").Value = 940
$ws.Range("T201").Value = 420

# Row 202
$ws.Range("D202").Value = 44258
$ws.Range("D202").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("K202").Value = 'June Pearl'
$ws.Range("L202").Value = 'Primera'
$ws.Range("M202").Value = 20
$ws.Range("N202").Value = 350000
$ws.Range("O202").Value = 360000
$ws.Range("P202").Value = 355000
$ws.Range("Q202").Value = '$/bins (420 kilos)'
$ws.Range("S202").Value = 845
$ws.Range("T202").Value = 420

# Row 203
$ws.Range("D203").Value = 44258
$ws.Range("D203").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("K203").Value = 'Venus'
$ws.Range("L203").Value = 'Especial'
$ws.Range("M203").Value = 16
$ws.Range("N203").Value = 360000
$ws.Range("O203").Value = 370000
$ws.Range("P203").Value = 365000
$ws.Range("Q203").Value = '$/bins (420 kilos)'
$ws.Range("S203").Value = 869
$ws.Range("T203").Value = 420

# Row 204
$ws.Range("D204").Value = 44595
$ws.Range("D204").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("K204").Value = 'Venus'
$ws.Range("L204").Value = 'Primera'
$ws.Range("M204").Value = 20
$ws.Range("N204").Value = 320000
$ws.Range("O204").Value = 330000
$ws.Range("P204").Value = 325000
$ws.Range("Q204").Value = '$/bins (420 kilos)'
$ws.Range("S204").Value = 774
$ws.Range("T204").Value = 420

# Row 205
$ws.Range("D205").Value = 44595
$ws.Range("D205").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("K205").Value = 'Venus'
$ws.Range("L205").Value = 'Segunda'
$ws.Range("M205").Value = 20
$ws.Range("N205").Value = 270000
$ws.Range("O205").Value = 280000
$ws.Range("P205").Value = 275000
$ws.Range("Q205").Value = '$/bins (420 kilos)'
$ws.Range("S205").Value = 655
$ws.Range("T205").Value = 420

# Row 206 (new)
$ws.Range("A206").Value = 2
$ws.Range("B206").Value = 'Comercializadora del Agro de Limarí'
$ws.Range("C206").Value = 'Coquimbo'
$ws.Range("D206").Value = 44552
$ws.Range("D206").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E206").Value = 4
$ws.Range("F206").Value = 'Fruta'
$ws.Range("G206").Value = 100103
$ws.Range("H206").Value = 'Frutos de hueso (carozo)'
$ws.Range("I206").Value = 100103006
$ws.Range("J206").Value = 'Nectarín'
$ws.Range("K206").Value = 'Big John'
$ws.Range("L206").Value = 'Especial'
$ws.Range("M206").Value = 24
$ws.Range("N206").Value = 485000
$ws.Range("O206").Value = 490000
$ws.Range("P206").Value = 487500
$ws.Range("Q206").Value = '$/bins (420 kilos)'
$ws.Range("R206").Value = 'Región de O''Higgins'
$ws.Range("S206").Value = 1161
$ws.Range("T206").Value = 420

# Row 207 (new)
$ws.Range("A207").Value = 2
$ws.Range("B207").Value = 'Comercializadora del Agro de Limarí'
$ws.Range("C207").Value = 'Coquimbo'
$ws.Range("D207").Value = 44552
$ws.Range("D207").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E207").Value = 4
$ws.Range("F207").Value = 'Fruta'
$ws.Range("G207").Value = 100103
$ws.Range("H207").Value = 'Frutos de hueso (carozo)'
$ws.Range("I207").Value = 100103006
$ws.Range("J207").Value = 'Nectarín'
$ws.Range("K207").Value = 'Big John'
$ws.Range("L207").Value = 'Primera'
$ws.Range("M207").Value = 20
$ws.Range("N207").Value = 435000
$ws.Range("O207").Value = 440000
$ws.Range("P207").Value = 437500
$ws.Range("Q207").Value = '$/bins (420 kilos)'
$ws.Range("R207").Value = 'Región de O''Higgins'
$ws.Range("S207").Value = 1042
$ws.Range("T207").Value = 420
